$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 145.4
$ws.Range("I33").Value = 130
$ws.Range("K33").Value = 130
$ws.Range("M33").Value = 99
# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 6290.0386
$ws.Range("I40").Value = 5271.9375
$ws.Range("K40").Value = 5271.9375
$ws.Range("M40").Value = -5096.9375
# row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 2227.5715
$ws.Range("J58").Value = 4504.25
$ws.Range("L58").Value = 13512.75
$ws.Range("N58").Value = -13812.75
# row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 3291.1667
$ws.Range("I106").Value = 3282.8
$ws.Range("K106").Value = 3282.8
$ws.Range("M106").Value = -2651.8
# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 14003.25
$ws.Range("I132").Value = 12621.471
$ws.Range("J132").Value = 21833.334
$ws.Range("K132").Value = 37864.413
$ws.Range("L132").Value = 65500.00199999999
$ws.Range("M132").Value = -35334.413
$ws.Range("N132").Value = -70560.00199999999
# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3012.7144
$ws.Range("I138").Value = 829
$ws.Range("J138").Value = 4905.2666
$ws.Range("K138").Value = 2487
$ws.Range("L138").Value = 14715.7998
$ws.Range("M138").Value = 2653
$ws.Range("N138").Value = -24995.7998

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 11 (Leve Item ID 3767)
$ws.Range("H11").Value = 983.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 983.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 983.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -1271.5
# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3370
$ws.Range("I32").Value = 3307.3076
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 3307.3076
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -3020.3076
$ws.Range("N32").Value = -5574
# row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 789.2857
$ws.Range("I97").Value = 464.53845
$ws.Range("K97").Value = 464.53845
$ws.Range("M97").Value = 31.46154999999999
# row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4313.5
$ws.Range("I132").Value = 4307.2964
$ws.Range("K132").Value = 12921.8892
$ws.Range("M132").Value = -10391.8892

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 2375
$ws.Range("I22").Value = 2375
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2375
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2202
$ws.Range("N22").ClearContents()
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 5537.5386
$ws.Range("I86").Value = 2317.6
$ws.Range("K86").Value = 2317.6
$ws.Range("M86").Value = -1194.6
# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 5537.5386
$ws.Range("I89").Value = 2317.6
$ws.Range("K89").Value = 11588
$ws.Range("M89").Value = -5972

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 2 (Leve Item ID 1820)
$ws.Range("H2").Value = 294.57144
$ws.Range("I2").Value = 191.25
$ws.Range("J2").Value = 432.33334
$ws.Range("K2").Value = 191.25
$ws.Range("L2").Value = 432.33334
$ws.Range("M2").Value = -78.25
$ws.Range("N2").Value = -658.33334
# row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1560.125
$ws.Range("I16").Value = 1411.8334
$ws.Range("K16").Value = 1411.8334
$ws.Range("M16").Value = -1124.8334
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4052.125
$ws.Range("I31").Value = 1952.2759
$ws.Range("K31").Value = 1952.2759
$ws.Range("M31").Value = -1657.2759
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4052.125
$ws.Range("I34").Value = 1952.2759
$ws.Range("K34").Value = 1952.2759
$ws.Range("M34").Value = -1750.2759
# row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 58417.25
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 58417.25
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1560.125
$ws.Range("I113").Value = 1411.8334
$ws.Range("K113").Value = 1411.8334
$ws.Range("M113").Value = 758.1666
# row 118 (Leve Item ID 26112)
$ws.Range("H118").Value = 65550
$ws.Range("J118").Value = 65550
$ws.Range("L118").Value = 65550
$ws.Range("N118").Value = -68864

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 62 (Leve Item ID 12845)
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# row 65 (Leve Item ID 12845)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 2307.6667
$ws.Range("J114").Value = 2769.8572
$ws.Range("L114").Value = 8309.571599999999
$ws.Range("N114").Value = -14817.5716

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 4132.222
$ws.Range("I70").Value = 4099.143
$ws.Range("J70").Value = 4248
$ws.Range("K70").Value = 4099.143
$ws.Range("L70").Value = 4248
$ws.Range("M70").Value = -3829.143
$ws.Range("N70").Value = -4788
# row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 4132.222
$ws.Range("I73").Value = 4099.143
$ws.Range("J73").Value = 4248
$ws.Range("K73").Value = 4099.143
$ws.Range("L73").Value = 4248
$ws.Range("M73").Value = -3163.143
$ws.Range("N73").Value = -6120
# row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 2212.7083
$ws.Range("I102").Value = 2091.5217
$ws.Range("K102").Value = 2091.5217
$ws.Range("M102").Value = -469.5216999999998
# row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 3133.4
$ws.Range("I122").Value = 3116.75
$ws.Range("K122").Value = 9350.25
$ws.Range("M122").Value = -6900.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 6706.5454
$ws.Range("J7").Value = 6607.75
$ws.Range("L7").Value = 6607.75
$ws.Range("N7").Value = -6831.75
# row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 3775.625
$ws.Range("I40").Value = 3775.625
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3775.625
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3639.625
$ws.Range("N40").ClearContents()
# row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 3572.4
$ws.Range("I61").Value = 2049
$ws.Range("K61").Value = 2049
$ws.Range("M61").Value = -1847
# row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
# row 101 (Leve Item ID 18549)
$ws.Range("H101").Value = 25750
$ws.Range("J101").Value = 25750
$ws.Range("L101").Value = 25750
$ws.Range("N101").Value = -32240
# row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 3572.4
$ws.Range("I113").Value = 2049
$ws.Range("K113").Value = 2049
$ws.Range("M113").Value = 121
# row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 6706.5454
$ws.Range("J126").Value = 6607.75
$ws.Range("L126").Value = 19823.25
$ws.Range("N126").Value = -24763.25
# row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 4738.8
$ws.Range("I132").Value = 4922.75
$ws.Range("J132").Value = 4003
$ws.Range("K132").Value = 14768.25
$ws.Range("L132").Value = 12009
$ws.Range("M132").Value = -12238.25
$ws.Range("N132").Value = -17069

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 10 (Leve Item ID 3313)
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
# row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 4586.316
$ws.Range("I126").Value = 2249
$ws.Range("K126").Value = 6747
$ws.Range("M126").Value = -4277

